$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 518, pushing existing rows 518:587 down to 519:588
$ws.Rows.Item(518).Insert()

# Populate the newly inserted row with the new data record
$ws.Cells.Item(518, 1).Value = 3
$ws.Cells.Item(518, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(518, 3).Value = "Coquimbo"
$ws.Cells.Item(518, 4).Value = 45131
$ws.Cells.Item(518, 5).Value = 5
$ws.Cells.Item(518, 6).Value = 100112009
$ws.Cells.Item(518, 7).Value = "Acelga"
$ws.Cells.Item(518, 8).Value = "Sin especificar"
$ws.Cells.Item(518, 9).Value = "Primera"
$ws.Cells.Item(518, 10).Value = 120
$ws.Cells.Item(518, 11).Value = 3000
$ws.Cells.Item(518, 12).Value = 3000
$ws.Cells.Item(518, 13).Value = 3000
$ws.Cells.Item(518, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(518, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(518, 16).Value = 500
$ws.Cells.Item(518, 17).Value = 6
$ws.Cells.Item(518, 18).Value = "Hortaliza"
